# Fix the formulas for the totals (and the underlying running-balance
# formulas that had drifted to the wrong source rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("Starting Balance" row) should carry forward from row 12, not row 9.
$ws.Range("D13").Formula = "=F12"
$ws.Range("E13").Formula = "=G12"

# Row 14 should carry forward from row 13 (which itself now correctly carries
# forward from row 12). Re-key the whole contiguous D14:D16 / E14:E16 run so
# that the running-balance formulas stay consistent down through row 16.
$ws.Range("D14:D16").Formula = "=F13"
$ws.Range("E14:E16").Formula = "=G13"

# Fix up the grand-total row: the totals need to also add back in the
# opening balance (D4/E4) that isn't otherwise included in the H4:J17 sums.
$ws.Range("H18").Formula = "=SUM(H4:H17)+D4"
$ws.Range("I18").Formula = "=SUM(I4:I17)+(E4-D4)"
$ws.Range("J18").Formula = "=SUM(J4:J17)+E4"

# Restore the active selection to match the authored state.
$ws.Range("J20").Select() | Out-Null
